$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The TC3 block (rows 22-26) and TC4 block (rows 29-33) swap their
# "Steps" / "Expected Results" content while keeping the TC3/TC4 labels
# (B22 and B29) in place.
#
# Before:
#   TC3 Steps (B26)            = "Chefe/Beneficiário Clica para realizar o empenho de uma diária."
#   TC3 Expected Results (D26) = "SYSTEM Apresenta a tela de Registrar Empenho"
#   TC4 Steps (B33)            = "Chefe/Beneficiário Clica para atribuir/desatribuir o registro a si mesmo."
#   TC4 Expected Results (D33) = "SYSTEM Atualiza a lista de registros de solicitações, onde o nome deverá constar o nome do usuário logado (que se atribuiu como responsável pelo empenho), no campo de atribuição (no caso de desatribuição, o nome deverá ser removido)."
#
# After:
#   TC3 Steps (B26)            = "Chefe/Beneficiário Clica para atribuir/desatribuir o registro a si mesmo."
#   TC3 Expected Results (D26) = "SYSTEM Atualiza a lista de registros de solicitações, onde o nome deverá constar o nome do usuário logado (que se atribuiu como responsável pelo empenho), no campo de atribuição (no caso de desatribuição, o nome deverá ser removido)."
#   TC4 Steps (B33)            = "Chefe/Beneficiário Clica para realizar o empenho de uma diária."
#   TC4 Expected Results (D33) = "SYSTEM Apresenta a tela de Registrar Empenho"

$ws.Range("B26").Value = "Chefe/Beneficiário Clica para atribuir/desatribuir o registro a si mesmo."
$ws.Range("D26").Value = "SYSTEM Atualiza a lista de registros de solicitações, onde o nome deverá constar o nome do usuário logado (que se atribuiu como responsável pelo empenho), no campo de atribuição (no caso de desatribuição, o nome deverá ser removido)."

$ws.Range("B33").Value = "Chefe/Beneficiário Clica para realizar o empenho de uma diária."
$ws.Range("D33").Value = "SYSTEM Apresenta a tela de Registrar Empenho"

$wb.Save()
